$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 11
# from 45170 (2023-09-01) to 45174 (2023-09-05)
$ws.Range("C2:C11").Value = 45174
